$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "55.778.80"
$ws.Range("E2").Value = "  +3.07%  "
$ws.Range("D3").Value = "2.497.13"
$ws.Range("E3").Value = "  +8.40%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "480.46"
$ws.Range("E5").Value = "  +7.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.93"
$ws.Range("E6").Value = "  +8.11%  "
$ws.Range("E7").Value = "  +0.40%  "
$ws.Range("E8").Value = "  +8.17%  "
$ws.Range("D9").Value = "2.496.96"
$ws.Range("E9").Value = "  +9.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0986"
$ws.Range("E10").Value = "  +7.00%  "
$ws.Range("E11").Value = "  +1.29%  "
$ws.Range("E12").Value = "  +5.12%  "
$ws.Range("E13").Value = "  +0.34%  "
$ws.Range("D14").Value = "2.932.53"
$ws.Range("E14").Value = "  +8.70%  "
$ws.Range("D15").Value = "55.793.95"
$ws.Range("E15").Value = "  +3.12%  "
$ws.Range("E16").Value = "  +13.67%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.45"
$ws.Range("E17").Value = "  +8.87%  "
$ws.Range("D18").Value = "2.503.83"
$ws.Range("E18").Value = "  +8.19%  "
$ws.Range("E19").Value = "  +7.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "320.92"
$ws.Range("E20").Value = "  +7.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.95"
$ws.Range("E21").Value = "  +5.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.997"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("E23").Value = "  +6.94%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "57.95"
$ws.Range("E24").Value = "  +4.13%  "
$ws.Range("E25").Value = "  +5.68%  "
$ws.Range("E26").Value = "  +3.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.403"
$ws.Range("E27").Value = "  +8.83%  "
$ws.Range("D28").Value = "2.609.29"
$ws.Range("E28").Value = "  +8.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.38"
$ws.Range("E29").Value = "  +8.49%  "
$ws.Range("D30").Value = "0.0₃0771"
$ws.Range("E30").Value = "  +8.82%  "
$ws.Range("E31").Value = "  +0.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "148.51"
$ws.Range("E32").Value = "  +1.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.11"
$ws.Range("E33").Value = "  +7.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.47"
$ws.Range("E34").Value = "  +9.31%  "
$ws.Range("E35").Value = "  +10.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.70"
$ws.Range("E36").Value = "  +2.43%  "
$ws.Range("E37").Value = "  +9.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.842"
$ws.Range("E38").Value = "  +0.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.33"
$ws.Range("E39").Value = "  +3.77%  "
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.62%  "
$ws.Range("B41").Value = "Mantle"
$ws.Range("C41").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.613"
$ws.Range("E41").Value = "  +18.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0549"
$ws.Range("E42").Value = "  +10.77%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.38"
$ws.Range("E43").Value = "  +7.23%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.31"
$ws.Range("E44").Value = "  +6.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.16"
$ws.Range("E45").Value = "  -1.33%  "
$ws.Range("D46").Value = "1.970.79"
$ws.Range("E46").Value = "  +2.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0904"
$ws.Range("E47").Value = "  +10.22%  "
$ws.Range("E48").Value = "  +7.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "250.16"
$ws.Range("E49").Value = "  +32.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.46"
$ws.Range("E50").Value = "  +11.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.47"
$ws.Range("E51").Value = "  +8.56%  "
